$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questions = @(
    @{
        r = 749
        B = 'FCM'
        C = 'CEFETMINAS'
        D = '2019'
        E = @'
É fundamental que o conhecimento tácito e o conhecimento explícito se relacionem em sinergia nas empresas e nos órgãos públicos.
Avalie as definições sobre a construção
( ) Socialização: o conhecimento tácito é repassado por meio de experiências.
( ) Exteriorização: o conhecimento tácito é traduzido em conceitos explícitos, por meio da utilização de metáforas, analogias e símbolos.
( ) Combinação: o conhecimento explícito é construído reunindo conhecimentos explícitos provenientes de várias fontes.de conhecimento.
( ) Internalização: é o final do ciclo, o conhecimento explícito, após ser internalizado, não passa novamente a ser um conhecimento tácito.
De acordo com as definições, a sequência correta é:
'@
        F = 'Conhecimentos Específicos'
        G = 'Gestão do Conhecimento e Tecnologia'
        H = 'Médio'
        I = 'ME'
        J = 'V, V, V, F'
        K = 'V, F, F, V'
        L = 'F, V, F, F'
        M = 'F, F, V, V'
        N = $null
        O = 'A'
        RTrailingBlank = $false
    },
    @{
        r = 750
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'O que é o conceito de bullwhip-effect (efeito chicote) em uma cadeia de suprimento?'
        F = 'Conhecimentos Específicos'
        G = 'Gestão da Cadeia de Suprimentos'
        H = 'Médio'
        I = 'ME'
        J = 'Um tipo de equipamento de transporte utilizado para movimentar mercadorias em longas distâncias.'
        K = 'A variação sazonal na demanda por produtos.'
        L = 'A amplificação da demanda, à medida que se move da ponta do consumidor de uma cadeia de suprimento para o fornecedor.'
        M = 'Um método de avaliação de fornecedores com base em suas classificações.'
        N = 'A velocidade com que os produtos são entregues aos clientes.'
        O = 'C'
        RTrailingBlank = $false
    },
    @{
        r = 751
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'As três dimensões da sustentabilidade, segundo o conceito de triple bottom line, são as seguintes:'
        F = 'Conhecimentos Específicos'
        G = 'Sustentabilidade'
        H = 'Médio'
        I = 'ME'
        J = 'econômica, política e ambiental'
        K = 'social, econômica e cultural'
        L = 'ambiental, econômica e social'
        M = 'cultural, ambiental e social'
        N = 'política, cultural e social.'
        O = 'C'
        RTrailingBlank = $false
    },
    @{
        r = 752
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'O seguinte critério é usado para avaliar a qualidade e a maturidade de um sistema organizacional na engenharia organizacional:'
        F = 'Conhecimentos Específicos'
        G = 'Engenharia Organizacional'
        H = 'Médio'
        I = 'ME'
        J = 'Critério EFQM'
        K = 'Critério CMMI'
        L = 'Critério ISO 9001'
        M = 'Critério PDCA'
        N = 'Critério COBIT'
        O = 'A'
        RTrailingBlank = $false
    },
    @{
        r = 753
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'Qual abordagem de engenharia organizacional se concentra na adequação entre contexto, estrutura e desempenho das empresas?'
        F = 'Conhecimentos Específicos'
        G = 'Engenharia Organizacional'
        H = 'Médio'
        I = 'ME'
        J = 'Teoria das Relações Humanas'
        K = 'Teoria da Contingência'
        L = 'Teoria institucional'
        M = 'Gestão por Objetivos (MBO)'
        N = 'Benchmarking'
        O = 'B'
        RTrailingBlank = $false
    },
    @{
        r = 754
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = @'
O grande volume de dados gerados por sensores, máquinas e equipamentos na Indústria 4.0 é chamado, em sistemas de informação, de Big Data. Para funcionar de forma eficiente, o Big Data tem alguns atributos chamados de Vs do Big Data.
Esses atributos são
'@
        F = 'Conhecimentos Específicos'
        G = 'Gestão do Conhecimento e Tecnologia'
        H = 'Médio'
        I = 'ME'
        J = 'Velocidade, Volume, Variedade, Veracidade, Valor'
        K = 'Velocidade, Volume, Variedade, Veracidade, Vigência'
        L = 'Velocidade, Volume, Variedade, Veracidade, Valor, Viabilidade'
        M = 'Velocidade, Volume, Variedade, Veracidade, Valor, Viabilidade, Vigência'
        N = 'Velocidade, Variedade, Veracidade, Valor, Viabilidade, Vigência'
        O = 'A'
        RTrailingBlank = $false
    },
    @{
        r = 755
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = @'
O termo Indústria 4.0 foi cunhado em 2011 na feira de Hanover, na Alemanha, para designar fábricas inteligentes onde os sistemas físicos e virtuais cooperam de forma global e flexível.
Essa colaboração é propiciada pelas seguintes tecnologias habilitadoras
'@
        F = 'Conhecimentos Específicos'
        G = 'Gestão do Conhecimento e Tecnologia'
        H = 'Médio'
        I = 'ME'
        J = 'Internet das Coisas ou Internet Of Things (Iot), Simulação Digital, Cyber-Security ou Segurança Digital, Integração de Sistemas'
        K = 'Robôs Autônomos, Internet das Coisas ou Internet Of Things (Iot), Cyber-Security ou Segurança Digital, Computação na Nuvem, Manufatura Aditiva, Big Data, Realidade Aumentada (“Augmented Reality”)'
        L = 'Robôs Autônomos, Internet das Coisas ou Internet Of Things (Iot), Simulação Digital, Cyber-Security ou Segurança Digital, Integração de Sistemas, Computação na Nuvem, Manufatura Aditiva, Big Data.'
        M = 'Robôs Autônomos, Internet das Coisas ou Internet Of Things (Iot), Simulação Digital, Integração de Sistemas, Computação na Nuvem, Manufatura Aditiva, Big Data, Realidade Aumentada (“Augmented Reality”).'
        N = 'Robôs Autônomos, Internet das Coisas ou Internet Of Things (Iot), Simulação Digital, Cyber-Security ou Segurança Digital, Integração de Sistemas, Computação na Nuvem, Manufatura Aditiva, Big Data, Realidade Aumentada (“Augmented Reality”).'
        O = 'E'
        RTrailingBlank = $false
    },
    @{
        r = 756
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'O bom gerenciamento da cadeia de suprimentos é fundamental para as empresas, por tratar-se da(o)'
        F = 'Conhecimentos Específicos'
        G = 'Gestão da Cadeia de Suprimentos'
        H = 'Médio'
        I = 'ME'
        J = 'estrutura organizacional que define as responsabilidades, as funções e as relações entre os diferentes agentes da logística empresarial.('
        K = 'rede de organizações envolvidas na produção, na transformação e na distribuição de um produto ou serviço desde os fornecedores até os clientes finais.('
        L = 'sistema de informação que integra e coordena os dados e os fluxos de informação entre os diversos elos da logística empresarial.('
        M = 'conjunto de processos que ligam as atividades de planejamento, de execução e de controle da logística empresarial em uma organização.('
        N = 'modelo de gestão que busca otimizar o desempenho da logística empresarial por meio da melhoria contínua dos processos e da qualidade'
        O = 'B'
        RTrailingBlank = $false
    },
    @{
        r = 757
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'Sobre a diferença entre logística inbound e outbound, tem-se que a logística inbound é a que se refere ao fluxo de bens e serviços'
        F = 'Conhecimentos Específicos'
        G = 'Logística'
        H = 'Médio'
        I = 'ME'
        J = 'que entram na organização, enquanto a logística outbound é a que se refere ao fluxo de bens e serviços que saem da organização.'
        K = 'que saem da organização, enquanto a logística outbound é a que se refere ao fluxo de bens e serviços que entram na organização.'
        L = 'entre a organização e seus fornecedores, enquanto a logística outbound é a que se refere ao fluxo de bens e serviços entre a organização e seus clientes.'
        M = 'entre a organização e seus clientes, enquanto a logística outbound é a que se refere ao fluxo de bens e serviços entre a organização e seus fornecedores.'
        N = 'dentro da organização, enquanto a logística outbound é a que se refere ao fluxo de bens e serviços fora da organização.'
        O = 'C'
        RTrailingBlank = $false
    },
    @{
        r = 758
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'Um Centro de Distribuição (CD) é uma'
        F = 'Conhecimentos Específicos'
        G = 'Logística'
        H = 'Médio'
        I = 'ME'
        J = 'rede de canais de distribuição que conecta os produtores aos consumidores finais.'
        K = 'área geográfica onde se concentram os clientes potenciais ou efetivos de um determinado produto ou serviço.'
        L = 'instalação física onde são armazenados temporariamente os produtos acabados antes de serem enviados aos clientes finais.'
        M = 'estratégia de marketing que visa aumentar a disponibilidade, a acessibilidade e a visibilidade dos produtos no mercado.'
        N = 'unidade operacional onde são realizadas as atividades de recebimento, de armazenagem, de separação, de embalagem e de expedição dos produtos ao longo da cadeia de suprimentos.'
        O = 'E'
        RTrailingBlank = $false
    },
    @{
        r = 759
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = @'
O sistema dutoviário de transportes no Brasil é aquele que utiliza tubulações para transportar produtos. Esse modal surgiu no país na década de 50 e tem algumas vantagens, como baixo custo operacional, alta capacidade de carga, menor emissão de poluentes, maior segurança e funcionamento contínuo. Por outro lado, também tem algumas desvantagens, como baixa velocidade, pouca flexibilidade de destinos e de produtos e infraestrutura precária.
Considerando-se o exposto, quais são os principais produtos transportados pelo modal dutoviário no Brasil?
'@
        F = 'Conhecimentos Específicos'
        G = 'Logística'
        H = 'Médio'
        I = 'ME'
        J = 'Óleo, carvão, cimento e água'
        K = 'Óleo, gás, minério e sal-gema'
        L = 'Óleo, minério, sal-gema e carvão'
        M = 'Gás, sal-gema, cimento e água'
        N = 'Minério, carvão, cimento e água.'
        O = 'B'
        RTrailingBlank = $false
    },
    @{
        r = 760
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = @'
Buscando atender à demanda de clientes, a manufatura utiliza três principais estratégias: antecipação da demanda, aumentando a capacidade produtiva antes que a demanda ocorra; seguimento da demanda, aguardando que a demanda ultrapasse a capacidade produtiva para aumentar a capacidade; e uma estratégia mista, combinando antecipação e seguimento da demanda.
Em indústrias de fluxo contínuo de produção, como siderúrgicas ou fábricas de celulose, a estratégia mais utilizada é a de seguir a demanda, pois nesse tipo de indústria
'@
        F = 'Conhecimentos Específicos'
        G = 'Gestão da Produção e Operações'
        H = 'Médio'
        I = 'ME'
        J = 'a complexidade do processo produtivo em fluxo contínuo exige maior tempo de planejamento para expandir.'
        K = 'a demanda tem uma grande elasticidade, e os cliente aguardam as expansões.'
        L = 'os tempos de realização da expansão são demasiadamente longos.'
        M = 'os valores financeiros necessários para a expansão são muito altos, tornando a antecipação inviável.'
        N = 'os fornecedores de máquinas e equipamentos para a indústria de fluxo contínuo demandam maior tempo de entrega.'
        O = 'D'
        RTrailingBlank = $false
    },
    @{
        r = 761
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'Qual é a diferença entre missão e visão no planejamento estratégico?'
        F = 'Conhecimentos Específicos'
        G = 'Engenharia Organizacional'
        H = 'Fácil'
        I = 'ME'
        J = 'Missão é a razão de ser da empresa, e visão é a situação desejada para o futuro.('
        K = 'Missão é a situação desejada para o futuro, e visão é a razão de ser da empresa.('
        L = 'Missão é o conjunto de valores que guiam as decisões da empresa, e visão é o diferencial competitivo que ela oferece.('
        M = 'Missão é o diferencial competitivo que a empresa oferece, e visão é o conjunto de valores que guiam as decisões da empresa.('
        N = 'Missão e visão são sinônimos no planejamento estratégico.'
        O = 'A'
        RTrailingBlank = $false
    },
    @{
        r = 762
        B = 'Cesgranrio'
        C = 'Transpetro'
        D = '2023'
        E = 'No planejamento estratégico, quais são as características das metas SMART?'
        F = 'Conhecimentos Específicos'
        G = 'Engenharia Organizacional'
        H = 'Médio'
        I = 'ME'
        J = 'Simples, Mensuráveis, Alcançáveis, Relevantes e Temporais'
        K = 'Específicas, Mensuráveis, Atribuíveis, Realistas e Temporais'
        L = 'Específicas, Mensuráveis, Alcançáveis, Relevantes e Temporais'
        M = 'Estratégicas, Mensuráveis, Ajustáveis, Realistas e Temporais'
        N = 'Estratégicas, Mensuráveis, Alcançáveis, Responsáveis e Temporais'
        O = 'C'
        RTrailingBlank = $true
    }
)
# Column letters C2..N2 map to spreadsheet columns B..N (2..14); O/P/Q are 15/16/17; R is 18.
$colIndex = @{ B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14; O=15 }

# R748 is already a "declared but blank" placeholder cell in the source sheet (t="inlineStr"
# with no text). Copying it elsewhere preserves that "present but empty" cell shape instead
# of Excel simply omitting the cell the way a plain `.Value = ""` on a brand new cell would.
$blankTemplate = $ws.Cells.Item(748, 18)

foreach ($q in $questions) {
    $r = $q.r

    $ws.Cells.Item($r, 1).Value = $r

    foreach ($col in @("B","C","D","E","F","G","H","I","J","K","L","M","N","O")) {
        $value = $q[$col]
        $target = $ws.Cells.Item($r, $colIndex[$col])
        if ($value -eq $null) {
            $blankTemplate.Copy($target)
        } else {
            $target.Value = $value
        }
    }

    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0

    if ($q.RTrailingBlank) {
        $blankTemplate.Copy($ws.Cells.Item($r, 18))
    }
}

# The old last row's trailing "R748" placeholder cell is dropped in the new layout (the
# blank marker moves to the very end of the sheet, on R762, instead).
$ws.Cells.Item(748, 18).Value = ""

# Re-fit the new rows' heights: Excel auto-expands row height when a multi-line value is
# assigned via COM, but the source rows all use the sheet's default height.
$ws.Rows("749:762").AutoFit()

Write-Output "Inserted rows 749-762 and cleared the old R748 placeholder."
